$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5933.849
$ws.Range("J17").Value = 5995.5
$ws.Range("L17").Value = 17986.5
$ws.Range("N17").Value = -18322.5
$ws.Range("H38").Value = 901.5833
$ws.Range("I38").Value = 81.90000000000001
$ws.Range("K38").Value = 245.7
$ws.Range("M38").Value = 126.3
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 1000
$ws.Range("K49").Value = 3000
$ws.Range("M49").Value = -2864
$ws.Range("H69").Value = 100502550
$ws.Range("I69").Value = 1668470
$ws.Range("K69").Value = 5005410
$ws.Range("M69").Value = -5004536
$ws.Range("H72").Value = 100502550
$ws.Range("I72").Value = 1668470
$ws.Range("K72").Value = 15016230
$ws.Range("M72").Value = -15011862
$ws.Range("H132").Value = 27030816
$ws.Range("I132").Value = 29415116
$ws.Range("K132").Value = 88245348
$ws.Range("M132").Value = -88242818
$ws.Range("H135").Value = 995.3570999999999
$ws.Range("J135").Value = 1742
$ws.Range("L135").Value = 15678
$ws.Range("N135").Value = -20748
$ws.Range("H137").Value = 83237.87
$ws.Range("I137").Value = 101040.945
$ws.Range("K137").Value = 303122.835
$ws.Range("M137").Value = -300572.835
$ws.Range("H138").Value = 2897.2078
$ws.Range("I138").Value = 1040.6538
$ws.Range("J138").Value = 3843.6863
$ws.Range("K138").Value = 3121.9614
$ws.Range("L138").Value = 11531.0589
$ws.Range("M138").Value = 2018.0386
$ws.Range("N138").Value = -21811.0589

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6323.269
$ws.Range("I32").Value = 3669.766
$ws.Range("J32").Value = 19093.25
$ws.Range("K32").Value = 3669.766
$ws.Range("L32").Value = 19093.25
$ws.Range("M32").Value = -3382.766
$ws.Range("N32").Value = -19667.25
$ws.Range("H61").Value = 7812.5
$ws.Range("I61").Value = 7812.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 7812.5
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -7600.5
$ws.Range("H74").Value = 86775.60000000001
$ws.Range("I74").Value = 31445.395
$ws.Range("J74").Value = 269365.3
$ws.Range("K74").Value = 31445.395
$ws.Range("L74").Value = 269365.3
$ws.Range("M74").Value = -30571.395
$ws.Range("N74").Value = -271113.3
$ws.Range("H77").Value = 86775.60000000001
$ws.Range("I77").Value = 31445.395
$ws.Range("J77").Value = 269365.3
$ws.Range("K77").Value = 157226.975
$ws.Range("L77").Value = 1346826.5
$ws.Range("M77").Value = -152858.975
$ws.Range("N77").Value = -1355562.5
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H122").Value = 908568.25
$ws.Range("I122").Value = 1997.5454
$ws.Range("J122").Value = 1739591.4
$ws.Range("K122").Value = 5992.6362
$ws.Range("L122").Value = 5218774.199999999
$ws.Range("M122").Value = -3542.6362
$ws.Range("N122").Value = -5223674.199999999
$ws.Range("H132").Value = 4650.5127
$ws.Range("I132").Value = 4588.9697
$ws.Range("J132").Value = 4989
$ws.Range("K132").Value = 13766.9091
$ws.Range("L132").Value = 14967
$ws.Range("M132").Value = -11236.9091
$ws.Range("N132").Value = -20027
$ws.Range("H136").Value = 7812.5
$ws.Range("I136").Value = 7812.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 23437.5
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -20887.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25048.49
$ws.Range("I31").Value = 7187.6665
$ws.Range("K31").Value = 7187.6665
$ws.Range("M31").Value = -6892.6665
$ws.Range("H34").Value = 25048.49
$ws.Range("I34").Value = 7187.6665
$ws.Range("K34").Value = 7187.6665
$ws.Range("M34").Value = -6985.6665
$ws.Range("H58").Value = 15747.25
$ws.Range("I58").Value = 26997
$ws.Range("K58").Value = 26997
$ws.Range("M58").Value = -26794
$ws.Range("H109").Value = 17996.334
$ws.Range("J109").Value = 17996.334
$ws.Range("L109").Value = 17996.334
$ws.Range("N109").Value = -20076.334
$ws.Range("H132").Value = 135150.33
$ws.Range("I132").Value = 93033.37
$ws.Range("K132").Value = 279100.11
$ws.Range("M132").Value = -276570.11
$ws.Range("H134").Value = 38589.387
$ws.Range("I134").Value = 64788.125
$ws.Range("K134").Value = 194364.375
$ws.Range("M134").Value = -191829.375
$ws.Range("H136").Value = 15747.25
$ws.Range("I136").Value = 26997
$ws.Range("K136").Value = 80991
$ws.Range("M136").Value = -78441

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5383202
$ws.Range("I4").Value = 8346096
$ws.Range("K4").Value = 25038288
$ws.Range("M4").Value = -25038176
$ws.Range("H12").Value = 63675.145
$ws.Range("I12").Value = 222387.25
$ws.Range("J12").Value = 190.3
$ws.Range("K12").Value = 667161.75
$ws.Range("L12").Value = 570.9000000000001
$ws.Range("M12").Value = -666988.75
$ws.Range("N12").Value = -916.9000000000001
$ws.Range("H40").Value = 53.526318
$ws.Range("I40").Value = 50.53846
$ws.Range("K40").Value = 202.15384
$ws.Range("M40").Value = -133.15384
$ws.Range("H62").Value = 2906
$ws.Range("I62").Value = 1812
$ws.Range("K62").Value = 5436
$ws.Range("M62").Value = -4750
$ws.Range("H63").Value = 13050
$ws.Range("I63").Value = 1850
$ws.Range("J63").Value = 24250
$ws.Range("K63").Value = 5550
$ws.Range("L63").Value = 72750
$ws.Range("M63").Value = -4801
$ws.Range("N63").Value = -74248
$ws.Range("H64").Value = 1500
$ws.Range("J64").Value = 1500
$ws.Range("L64").Value = 4500
$ws.Range("N64").Value = -5040
$ws.Range("H65").Value = 2906
$ws.Range("I65").Value = 1812
$ws.Range("K65").Value = 16308
$ws.Range("M65").Value = -12876
$ws.Range("H66").Value = 13050
$ws.Range("I66").Value = 1850
$ws.Range("J66").Value = 24250
$ws.Range("K66").Value = 16650
$ws.Range("L66").Value = 218250
$ws.Range("M66").Value = -12906
$ws.Range("N66").Value = -225738
$ws.Range("H67").Value = 1500
$ws.Range("J67").Value = 1500
$ws.Range("L67").Value = 4500
$ws.Range("N67").Value = -6372
$ws.Range("H69").Value = 2222.4
$ws.Range("I69").Value = 1056
$ws.Range("K69").Value = 3168
$ws.Range("M69").Value = -2357
$ws.Range("H72").Value = 2222.4
$ws.Range("I72").Value = 1056
$ws.Range("K72").Value = 9504
$ws.Range("M72").Value = -5448
$ws.Range("H131").Value = 19846426
$ws.Range("I131").Value = 9262012
$ws.Range("J131").Value = 27784736
$ws.Range("K131").Value = 27786036
$ws.Range("L131").Value = 83354208
$ws.Range("M131").Value = -27780996
$ws.Range("N131").Value = -83364288

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 14962.5
$ws.Range("I132").Value = 10626.25
$ws.Range("J132").Value = 27971.25
$ws.Range("K132").Value = 31878.75
$ws.Range("L132").Value = 83913.75
$ws.Range("M132").Value = -29348.75
$ws.Range("N132").Value = -88973.75
$ws.Range("H141").Value = 34999.75
$ws.Range("J141").Value = 34999.75
$ws.Range("L141").Value = 34999.75
$ws.Range("N141").Value = -45359.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6292
$ws.Range("J7").Value = 9506.143
$ws.Range("L7").Value = 9506.143
$ws.Range("N7").Value = -9730.143
$ws.Range("H68").Value = 2862
$ws.Range("I68").Value = 2862
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2862
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -2113
$ws.Range("H71").Value = 2862
$ws.Range("I71").Value = 2862
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 14310
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -10566
$ws.Range("H82").Value = 5559346
$ws.Range("I82").Value = 7939066
$ws.Range("K82").Value = 7939066
$ws.Range("M82").Value = -7938705
$ws.Range("H85").Value = 5559346
$ws.Range("I85").Value = 7939066
$ws.Range("K85").Value = 7939066
$ws.Range("M85").Value = -7937818
$ws.Range("H126").Value = 6292
$ws.Range("J126").Value = 9506.143
$ws.Range("L126").Value = 28518.429
$ws.Range("N126").Value = -33458.429
$ws.Range("H133").Value = 98121.28999999999
$ws.Range("J133").Value = 98121.28999999999
$ws.Range("L133").Value = 98121.28999999999
$ws.Range("N133").Value = -103181.29
$ws.Range("H136").Value = 84534.44
$ws.Range("I136").Value = 156181.23
$ws.Range("K136").Value = 468543.6900000001
$ws.Range("M136").Value = -465993.6900000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 17996.715
$ws.Range("J31").Value = 21195.4
$ws.Range("L31").Value = 21195.4
$ws.Range("N31").Value = -21891.4
$ws.Range("H136").Value = 5059.204
$ws.Range("I136").Value = 6094.2666
$ws.Range("J136").Value = 3424.8948
$ws.Range("K136").Value = 18282.7998
$ws.Range("L136").Value = 10274.6844
$ws.Range("M136").Value = -15732.7998
$ws.Range("N136").Value = -15374.6844
